{"js": "// Append the \"Allegati\" (attachments) block to the end of the document body:\n//   <empty paragraph>\n//   \"Allegati\"               (bold)\n//   \"{#attachments}\"\n//   \"{image}\"\n//   \"{/attachments}\"\nconst body = context.document.body;\n\n// Queue up all five new paragraphs first (in document order) - the new\n// paragraph mark otherwise picks up whatever direct formatting was just\n// applied to the previous paragraph, so we only touch \"bold\" after every\n// paragraph already exists.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nconst heading = body.insertParagraph(\"Allegati\", Word.InsertLocation.end);\nbody.insertParagraph(\"{#attachments}\", Word.InsertLocation.end);\nbody.insertParagraph(\"{image}\", Word.InsertLocation.end);\nbody.insertParagraph(\"{/attachments}\", Word.InsertLocation.end);\n\n// Now make the heading paragraph bold.\nheading.font.bold = true;\n\nawait context.sync();\n", "ps1": "# Append the \"Allegati\" (attachments) block to the end of the document:\n#   <empty paragraph>\n#   \"Allegati\"               (bold)\n#   \"{#attachments}\"\n#   \"{image}\"\n#   \"{/attachments}\"\n$d = $word.ActiveDocument\n\n# Add all five paragraphs first (Paragraphs.Add() appends at the end of the\n# story) before touching any direct character formatting - Word carries the\n# last-applied direct formatting forward onto the next freshly-added\n# paragraph mark, so only the \"Allegati\" run should be made bold once every\n# paragraph already exists.\n$pBlank = $d.Paragraphs.Add()\n$pBlank.Range.Text = \"\"\n\n$pHeading = $d.Paragraphs.Add()\n$pHeading.Range.Text = \"Allegati\"\n\n$pOpen = $d.Paragraphs.Add()\n$pOpen.Range.Text = \"{#attachments}\"\n\n$pImage = $d.Paragraphs.Add()\n$pImage.Range.Text = \"{image}\"\n\n$pClose = $d.Paragraphs.Add()\n$pClose.Range.Text = \"{/attachments}\"\n\n# Bold only the visible \"Allegati\" text, not the trailing paragraph mark\n# (Range.End - 1), so the paragraph-mark formatting stays plain.\n$headingText = $d.Range($pHeading.Range.Start, $pHeading.Range.End - 1)\n$headingText.Bold = 1\n"}
